$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4856.75
$ws.Range("I51").Value = 3384.4285
$ws.Range("J51").Value = 6918
$ws.Range("K51").Value = 3384.4285
$ws.Range("L51").Value = 6918
$ws.Range("M51").Value = -2900.4285
$ws.Range("N51").Value = -7886
$ws.Range("H123").Value = 41811.25
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 41811.25
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 41811.25
$ws.Range("N123").Value = -51611.25
$ws.Range("H128").Value = 41871.25
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 41871.25
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 41871.25
$ws.Range("N128").Value = -51831.25
$ws.Range("H131").Value = 3878
$ws.Range("I131").Value = 3377.5
$ws.Range("J131").Value = 4879
$ws.Range("K131").Value = 10132.5
$ws.Range("L131").Value = 14637
$ws.Range("M131").Value = -5092.5
$ws.Range("N131").Value = -24717
$ws.Range("H134").Value = 52201.4
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 52201.4
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 52201.4
$ws.Range("N134").Value = -62341.4
$ws.Range("H137").Value = 1289059.1
$ws.Range("I137").Value = 1489546.5
$ws.Range("J137").Value = 5939.8
$ws.Range("K137").Value = 4468639.5
$ws.Range("L137").Value = 17819.4
$ws.Range("M137").Value = -4466089.5
$ws.Range("N137").Value = -22919.4
$ws.Range("H140").Value = 46711.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 46711.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 46711.332
$ws.Range("N140").Value = -57071.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H109").Value = 44500
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 44500
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 44500
$ws.Range("N109").Value = -47274
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 4965.16
$ws.Range("I122").Value = 4055.75
$ws.Range("J122").Value = 8602.799999999999
$ws.Range("K122").Value = 12167.25
$ws.Range("L122").Value = 25808.4
$ws.Range("M122").Value = -9717.25
$ws.Range("N122").Value = -30708.4
$ws.Range("H130").Value = 43928.5
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 43928.5
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 43928.5
$ws.Range("N130").Value = -53968.5
$ws.Range("H132").Value = 2866.6667
$ws.Range("I132").Value = 1448
$ws.Range("J132").Value = 3576
$ws.Range("K132").Value = 4344
$ws.Range("L132").Value = 10728
$ws.Range("M132").Value = -1814
$ws.Range("N132").Value = -15788
$ws.Range("H137").Value = 43753.332
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 43753.332
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 43753.332
$ws.Range("N137").Value = -53953.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H137").Value = 49537
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 49537
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 49537
$ws.Range("N137").Value = -59737

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9260910
$ws.Range("I16").Value = 15874559
$ws.Range("J16").Value = 1802.6
$ws.Range("K16").Value = 15874559
$ws.Range("L16").Value = 1802.6
$ws.Range("M16").Value = -15874272
$ws.Range("N16").Value = -2376.6
$ws.Range("H31").Value = 6263.5137
$ws.Range("I31").Value = 2735.3914
$ws.Range("J31").Value = 12059.714
$ws.Range("K31").Value = 2735.3914
$ws.Range("L31").Value = 12059.714
$ws.Range("M31").Value = -2440.3914
$ws.Range("N31").Value = -12649.714
$ws.Range("H34").Value = 6263.5137
$ws.Range("I34").Value = 2735.3914
$ws.Range("J34").Value = 12059.714
$ws.Range("K34").Value = 2735.3914
$ws.Range("L34").Value = 12059.714
$ws.Range("M34").Value = -2533.3914
$ws.Range("N34").Value = -12463.714
$ws.Range("H113").Value = 9260910
$ws.Range("I113").Value = 15874559
$ws.Range("J113").Value = 1802.6
$ws.Range("K113").Value = 15874559
$ws.Range("L113").Value = 1802.6
$ws.Range("M113").Value = -15872389
$ws.Range("N113").Value = -6142.6
$ws.Range("H122").Value = 2984.7778
$ws.Range("I122").Value = 1916.3334
$ws.Range("J122").Value = 3519
$ws.Range("K122").Value = 5749.0002
$ws.Range("L122").Value = 10557
$ws.Range("M122").Value = -3299.0002
$ws.Range("N122").Value = -15457
$ws.Range("H138").Value = 40758.332
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 40758.332
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 40758.332
$ws.Range("N138").Value = -51038.332
$ws.Range("H140").Value = 105568.89
$ws.Range("I140").Value = 10000
$ws.Range("J140").Value = 117515
$ws.Range("K140").Value = 10000
$ws.Range("L140").Value = 117515
$ws.Range("M140").Value = -4820
$ws.Range("N140").Value = -127875
$ws.Range("H141").Value = 32031.818
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 32031.818
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 32031.818
$ws.Range("N141").Value = -42391.818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 7536.2856
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 8625.666999999999
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 25877.001
$ws.Range("M62").Value = -2314
$ws.Range("N62").Value = -27249.001
$ws.Range("H63").Value = 4212.5
$ws.Range("I63").Value = 3511.111
$ws.Range("J63").Value = 4786.364
$ws.Range("K63").Value = 10533.333
$ws.Range("L63").Value = 14359.092
$ws.Range("M63").Value = -9784.332999999999
$ws.Range("N63").Value = -15857.092
$ws.Range("H64").Value = 1807.6154
$ws.Range("I64").Value = 833
$ws.Range("J64").Value = 2100
$ws.Range("K64").Value = 2499
$ws.Range("L64").Value = 6300
$ws.Range("M64").Value = -2229
$ws.Range("N64").Value = -6840
$ws.Range("H65").Value = 7536.2856
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 8625.666999999999
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 77631.003
$ws.Range("M65").Value = -5568
$ws.Range("N65").Value = -84495.003
$ws.Range("H66").Value = 4212.5
$ws.Range("I66").Value = 3511.111
$ws.Range("J66").Value = 4786.364
$ws.Range("K66").Value = 31599.999
$ws.Range("L66").Value = 43077.276
$ws.Range("M66").Value = -27855.999
$ws.Range("N66").Value = -50565.276
$ws.Range("H67").Value = 1807.6154
$ws.Range("I67").Value = 833
$ws.Range("J67").Value = 2100
$ws.Range("K67").Value = 2499
$ws.Range("L67").Value = 6300
$ws.Range("M67").Value = -1563
$ws.Range("N67").Value = -8172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9334938
$ws.Range("I11").Value = 16666742
$ws.Range("J11").Value = 2003134.1
$ws.Range("K11").Value = 16666742
$ws.Range("L11").Value = 2003134.1
$ws.Range("M11").Value = -16666603
$ws.Range("N11").Value = -2003412.1
$ws.Range("H33").Value = 4152.8335
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 4383.4
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 4383.4
$ws.Range("M33").Value = -2748
$ws.Range("N33").Value = -4887.4
$ws.Range("H42").Value = 27356
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 27356
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 27356
$ws.Range("N42").Value = -28326
$ws.Range("H115").Value = 27356
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 27356
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 27356
$ws.Range("N115").Value = -29706
$ws.Range("H122").Value = 5903.3335
$ws.Range("I122").Value = 4084
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 12252
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -9802
$ws.Range("N122").Value = -49900
$ws.Range("H126").Value = 3271.4895
$ws.Range("I126").Value = 2819.25
$ws.Range("J126").Value = 4990
$ws.Range("K126").Value = 8457.75
$ws.Range("L126").Value = 14970
$ws.Range("M126").Value = -5987.75
$ws.Range("N126").Value = -19910
$ws.Range("H132").Value = 3790
$ws.Range("I132").Value = 2369.2727
$ws.Range("J132").Value = 8999.333000000001
$ws.Range("K132").Value = 7107.8181
$ws.Range("L132").Value = 26997.999
$ws.Range("M132").Value = -4577.8181
$ws.Range("N132").Value = -32057.999
$ws.Range("H137").Value = 40193.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 40193.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 40193.5
$ws.Range("N137").Value = -50393.5
$ws.Range("H141").Value = 40650
$ws.Range("I141").Value = 36000
$ws.Range("J141").Value = 42975
$ws.Range("K141").Value = 36000
$ws.Range("L141").Value = 42975
$ws.Range("M141").Value = -30820
$ws.Range("N141").Value = -53335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11633.333
$ws.Range("I40").Value = 9950
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 9950
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -9814
$ws.Range("N40").Value = -15272
$ws.Range("H46").Value = 1542.5
$ws.Range("I46").Value = 926.6667
$ws.Range("J46").Value = 2158.3333
$ws.Range("K46").Value = 926.6667
$ws.Range("L46").Value = 2158.3333
$ws.Range("M46").Value = -738.6667
$ws.Range("N46").Value = -2534.3333
$ws.Range("H98").Value = 35000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 35000
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -40990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 41213.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 41213.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 41213.75
$ws.Range("N135").Value = -51353.75
$ws.Range("H139").Value = 39487.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 39487.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 39487.5
$ws.Range("N139").Value = -49767.5
$ws.Range("H140").Value = 34836.4
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 34836.4
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 34836.4
$ws.Range("N140").Value = -45196.4
$ws.Range("H141").Value = 37256
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 37256
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 37256
$ws.Range("N141").Value = -47616
